$d = $word.ActiveDocument

# --- Change 1: cover letter opening paragraph ---
# "I and my co-authors have discussed ... manuscript in response to each."
# becomes "My co-authors and I have discussed ... manuscript in response to each."
# (the _GoBack bookmark sits between "the " and "manuscript" and must be preserved)
$old1a = "I and my co-authors have discussed all reviewer comments, and have modified the text of the "
$new1a = "My co-authors and I "
$f1a = $d.Content.Find.Execute($old1a, $true, $false, $false, $false, $false, $true, 1, $false, $new1a, 2)
if (-not $f1a) { throw "change 1a: text not found" }

$old1b = "manuscript in response to each."
$new1b = "have discussed all reviewer comments, and have modified the text of the manuscript in response to each."
$f1b = $d.Content.Find.Execute($old1b, $true, $false, $false, $false, $false, $true, 1, $false, $new1b, 2)
if (-not $f1b) { throw "change 1b: text not found" }

# --- Change 2: delete the blank paragraph after the alpha-significance sentence, and
#     append the new sentence about the leave-one-out cross-validation to that paragraph ---
$old2 = "These differences are not statistically significant, because the number of tokens we were able to acquire in this experiment is not large enough to make this difference significant."
$r2 = $d.Content
$f2 = $r2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $f2) { throw "change 2: anchor text not found" }
$p2index = $r2.Paragraphs.First.Index
$blankPara = $d.Paragraphs($p2index + 1)
if ($blankPara.Range.Text.Length -gt 1) { throw "change 2: paragraph after anchor was not blank as expected" }
$blankPara.Range.Delete()

$add2 = "  The following text has been added: “The constant `$\alpha=0.29`$ was chosen as the average of the values selected in all folds of a leave-one-out cross-validation.”"
$r2b = $d.Content
$r2b.Find.Execute($old2, $false) | Out-Null
$r2b.Collapse(0)
$r2b.InsertAfter($add2)

# --- Change 3: no textual change, the quoted sentence about GMM/NN initialization is
#     simply re-typed as a single run in the source revision; re-apply identical text via
#     Find/Replace so the runs collapse the same way ---
$old3 = "GMM parameters were initialized using a monophone system trained on the same 40 minutes, NN parameters were initialized using a restricted Boltzmann machine trained on five hours of unlabeled audio in the same language."
$f3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)
if (-not $f3) { throw "change 3: text not found" }

# --- Change 4: same as Change 3, but for the Table III sentence ---
$old4 = "This sentence has been made more precise: “Table III showed that PTs computed using a text-based phone bigram language model achieve PER in the range 50.45-70.88%, depending on the language.”"
$f4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2)
if (-not $f4) { throw "change 4: text not found" }

Write-Host "All four changes applied successfully."
